$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 6449.8335
$ws.Range("I58").Value = 258
$ws.Range("J58").Value = 12641.667
$ws.Range("K58").Value = 774
$ws.Range("L58").Value = 37925.001
$ws.Range("M58").Value = -624
$ws.Range("N58").Value = -38225.001

$ws.Range("H137").Value = 1426.5957
$ws.Range("I137").Value = 1217.081
$ws.Range("J137").Value = 2201.8
$ws.Range("K137").Value = 3651.242999999999
$ws.Range("L137").Value = 6605.400000000001
$ws.Range("M137").Value = -1101.242999999999
$ws.Range("N137").Value = -11705.4

$ws.Range("H138").Value = 1776.85
$ws.Range("I138").Value = 1474.5483
$ws.Range("J138").Value = 2100
$ws.Range("K138").Value = 4423.644899999999
$ws.Range("L138").Value = 6300
$ws.Range("M138").Value = 716.3551000000007
$ws.Range("N138").Value = -16580

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 739.9535
$ws.Range("I74").Value = 654.4211
$ws.Range("J74").Value = 1390
$ws.Range("K74").Value = 654.4211
$ws.Range("L74").Value = 1390
$ws.Range("M74").Value = 219.5789
$ws.Range("N74").Value = -3138

$ws.Range("H77").Value = 739.9535
$ws.Range("I77").Value = 654.4211
$ws.Range("J77").Value = 1390
$ws.Range("K77").Value = 3272.1055
$ws.Range("L77").Value = 6950
$ws.Range("M77").Value = 1095.8945
$ws.Range("N77").Value = -15686

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4275668
$ws.Range("I31").Value = 1668.8182
$ws.Range("J31").Value = 9806726
$ws.Range("K31").Value = 1668.8182
$ws.Range("L31").Value = 9806726
$ws.Range("M31").Value = -1373.8182
$ws.Range("N31").Value = -9807316

$ws.Range("H34").Value = 4275668
$ws.Range("I34").Value = 1668.8182
$ws.Range("J34").Value = 9806726
$ws.Range("K34").Value = 1668.8182
$ws.Range("L34").Value = 9806726
$ws.Range("M34").Value = -1466.8182
$ws.Range("N34").Value = -9807130

$ws.Range("H58").Value = 7752782
$ws.Range("I58").Value = 810.8857400000001
$ws.Range("J58").Value = 41667656
$ws.Range("K58").Value = 810.8857400000001
$ws.Range("L58").Value = 41667656
$ws.Range("M58").Value = -607.8857400000001
$ws.Range("N58").Value = -41668062

$ws.Range("H136").Value = 7752782
$ws.Range("I136").Value = 810.8857400000001
$ws.Range("J136").Value = 41667656
$ws.Range("K136").Value = 2432.65722
$ws.Range("L136").Value = 125002968
$ws.Range("M136").Value = 117.3427799999999
$ws.Range("N136").Value = -125008068

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1187.52
$ws.Range("J5").Value = 1745.7333
$ws.Range("L5").Value = 5237.199900000001
$ws.Range("N5").Value = -5461.199900000001

$ws.Range("H63").Value = 7957.143
$ws.Range("I63").Value = 5425
$ws.Range("J63").Value = 11333.333
$ws.Range("K63").Value = 16275
$ws.Range("L63").Value = 33999.999
$ws.Range("M63").Value = -15526
$ws.Range("N63").Value = -35497.999

$ws.Range("H64").Value = 1250
$ws.Range("I64").Value = 1250
$ws.Range("K64").Value = 3750
$ws.Range("M64").Value = -3480

$ws.Range("H66").Value = 7957.143
$ws.Range("I66").Value = 5425
$ws.Range("J66").Value = 11333.333
$ws.Range("K66").Value = 48825
$ws.Range("L66").Value = 101999.997
$ws.Range("M66").Value = -45081
$ws.Range("N66").Value = -109487.997

$ws.Range("H67").Value = 1250
$ws.Range("I67").Value = 1250
$ws.Range("K67").Value = 3750
$ws.Range("M67").Value = -2814

$ws.Range("H68").Value = 308.7143
$ws.Range("I68").Value = 308.7143
$ws.Range("K68").Value = 926.1428999999999
$ws.Range("M68").Value = -115.1428999999999

$ws.Range("H71").Value = 308.7143
$ws.Range("I71").Value = 308.7143
$ws.Range("K71").Value = 2778.4287
$ws.Range("M71").Value = 1277.5713

$ws.Range("H94").Value = 9350.666999999999
$ws.Range("I94").Value = 3100
$ws.Range("J94").Value = 11434.223
$ws.Range("K94").Value = 9300
$ws.Range("L94").Value = 34302.669
$ws.Range("M94").Value = -8624
$ws.Range("N94").Value = -35654.669

$ws.Range("H103").Value = 6954.125
$ws.Range("I103").Value = 312.5
$ws.Range("J103").Value = 9168
$ws.Range("K103").Value = 937.5
$ws.Range("L103").Value = 27504
$ws.Range("M103").Value = -58.5
$ws.Range("N103").Value = -29262

$ws.Range("H108").Value = 180.75
$ws.Range("I108").Value = 180.75
$ws.Range("K108").Value = 542.25
$ws.Range("M108").Value = 2337.75

$ws.Range("H114").Value = 1168.9231
$ws.Range("I114").Value = 1081.5
$ws.Range("J114").Value = 1270.9166
$ws.Range("K114").Value = 3244.5
$ws.Range("L114").Value = 3812.7498
$ws.Range("M114").Value = 9.5
$ws.Range("N114").Value = -10320.7498

$ws.Range("H129").Value = 1470.4
$ws.Range("I129").Value = 920
$ws.Range("J129").Value = 1706.2858
$ws.Range("K129").Value = 2760
$ws.Range("L129").Value = 5118.857400000001
$ws.Range("M129").Value = 2240
$ws.Range("N129").Value = -15118.8574

$ws.Range("H134").Value = 1114.9546
$ws.Range("I134").Value = 752.4375
$ws.Range("K134").Value = 2257.3125
$ws.Range("M134").Value = 2812.6875

$ws.Range("H135").Value = 1187.52
$ws.Range("J135").Value = 1745.7333
$ws.Range("L135").Value = 15711.5997
$ws.Range("N135").Value = -20781.5997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2576.6667
$ws.Range("I126").Value = 2912.5
$ws.Range("J126").Value = 2192.8572
$ws.Range("K126").Value = 8737.5
$ws.Range("L126").Value = 6578.571599999999
$ws.Range("M126").Value = -6267.5
$ws.Range("N126").Value = -11518.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H88").Value = 44490
$ws.Range("J88").Value = 44490
$ws.Range("L88").Value = 44490
$ws.Range("N88").Value = -45346

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H91").Value = 44490
$ws.Range("J91").Value = 44490
$ws.Range("L91").Value = 44490
$ws.Range("N91").Value = -47454

$ws.Range("H132").Value = 5150.6875
$ws.Range("I132").Value = 6956.9033
$ws.Range("J132").Value = 1857
$ws.Range("K132").Value = 20870.7099
$ws.Range("L132").Value = 5571
$ws.Range("M132").Value = -18340.7099
$ws.Range("N132").Value = -10631

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1120.2709
$ws.Range("I132").Value = 1103.5897
$ws.Range("J132").Value = 1192.5555
$ws.Range("K132").Value = 3310.7691
$ws.Range("L132").Value = 3577.6665
$ws.Range("M132").Value = -780.7691
$ws.Range("N132").Value = -8637.666499999999

$ws.Range("H136").Value = 2504.5151
$ws.Range("I136").Value = 2713.9124
$ws.Range("J136").Value = 1178.3334
$ws.Range("K136").Value = 8141.7372
$ws.Range("L136").Value = 3535.0002
$ws.Range("M136").Value = -5591.7372
$ws.Range("N136").Value = -8635.0002
